$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @(1.465341735336888, 17.30416868637906, 1.019631604459525, -0.2947164323800291, 1.22200678623364, 1.207222491517208, 8.792777508482793, 9)
    3 = @(2.019693578208169, 9.195957851124549, 4.270133340904041, 0.9703013987511284, 1.085512506669432, 0.5297969353132786, 0.5297969353132786, 0)
    4 = @(1.475273418925345, 17.7246790299572, 0.9864479315263909, -0.3031383737517832, 1.187351402757965, 1.245601521635759, 8.75439847836424, 9)
    5 = @(2.026162321730315, 9.413694306879455, 4.273531832562288, 0.9636063076605779, 1.05410889547806, 0.5609770710404465, 0.5609770710404464, 0)
    6 = @(1.563603440470801, 18.29923042100498, 0.9637845178813943, -0.3198279764599366, 1.234721593481213, 1.279764599366229, 8.72023540063377, 9)
    7 = @(2.204438451590914, 10.67865166581983, 4.385752075066329, 1.019737640802872, 1.108306638962272, 0.649307977901109, 0.649307977901109, 0)
    8 = @(1.490718124084025, 17.66927210552027, 0.9201270151441133, -0.3094772838299951, 1.140693698094773, 1.256228627259404, 8.743771372740596, 9)
    9 = @(2.075954932242621, 9.857303339678452, 4.172053473022855, 0.9699733740136783, 1.03882138644175, 0.5894703072243467, 0.5894703072243467, 0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

Write-Output "Done updating cells"